# Update absenteeism data rows 2-11 with new values as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 39879
$ws.Range("B2").Value = "Juliana Silveira"
$ws.Range("C2").Value = "Engenharia"
$ws.Range("D2").Value = "Viagem de negócios"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 45094
$ws.Range("G2").Value = 4865.38

# Row 3
$ws.Range("A3").Value = 98497
$ws.Range("B3").Value = "Diogo Azevedo"
$ws.Range("C3").Value = "Atendimento ao Cliente"
$ws.Range("D3").Value = "Doença"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 45097
$ws.Range("G3").Value = 7556.64

# Row 4
$ws.Range("A4").Value = 4140
$ws.Range("B4").Value = "Larissa Lopes"
$ws.Range("C4").Value = "TI"
$ws.Range("D4").Value = "Outros"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 45093
$ws.Range("G4").Value = 4145.37

# Row 5
$ws.Range("A5").Value = 97887
$ws.Range("B5").Value = "Marcelo Vieira"
$ws.Range("C5").Value = "Marketing"
$ws.Range("D5").Value = "Consulta médica"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45080
$ws.Range("G5").Value = 4106.09

# Row 6
$ws.Range("A6").Value = 39737
$ws.Range("B6").Value = "Maria Sophia Santos"
$ws.Range("C6").Value = "Financeiro"
$ws.Range("D6").Value = "Doença"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 45090
$ws.Range("G6").Value = 9556.09

# Row 7
$ws.Range("A7").Value = 39642
$ws.Range("B7").Value = "Thiago Almeida"
$ws.Range("C7").Value = "Engenharia"
$ws.Range("D7").Value = "Consulta médica"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 45087
$ws.Range("G7").Value = 12045.56

# Row 8
$ws.Range("A8").Value = 54622
$ws.Range("B8").Value = "Dr. Guilherme Moura"
$ws.Range("C8").Value = "Jurídico"
$ws.Range("D8").Value = "Viagem de negócios"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 45080
$ws.Range("G8").Value = 3836.89

# Row 9
$ws.Range("A9").Value = 55082
$ws.Range("B9").Value = "Caroline Vieira"
$ws.Range("C9").Value = "Financeiro"
$ws.Range("D9").Value = "Problemas pessoais"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 45094
$ws.Range("G9").Value = 4531.74

# Row 10
$ws.Range("A10").Value = 87129
$ws.Range("B10").Value = "Luiz Fernando Rodrigues"
$ws.Range("C10").Value = "Financeiro"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 45091
$ws.Range("G10").Value = 6825.79

# Row 11
$ws.Range("A11").Value = 34256
$ws.Range("B11").Value = "Marcela Souza"
$ws.Range("C11").Value = "Financeiro"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 45079
$ws.Range("G11").Value = 2731.52
